$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix comma-as-separator typos in contractor/company names (E column) ---
$ws.Range("E55").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E66").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E114").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E56").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F56").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E57").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E115").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E69").Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("F69").Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA"

# --- Fix Importe (H column) formatting: "1.234,56" (es-AR) -> "1234.56" (plain decimal point, text) ---
$importeRange = $ws.Range("H2:H156")
$importeRange.NumberFormat = "@"
$ws.Range("H2").Value = "342.00"
$ws.Range("H3").Value = "7497.01"
$ws.Range("H4").Value = "9900.00"
$ws.Range("H5").Value = "123.00"
$ws.Range("H6").Value = "19.80"
$ws.Range("H7").Value = "7800.00"
$ws.Range("H8").Value = "78199.96"
$ws.Range("H9").Value = "34201.67"
$ws.Range("H10").Value = "5319.50"
$ws.Range("H11").Value = "680.00"
$ws.Range("H12").Value = "4407.64"
$ws.Range("H13").Value = "4655.00"
$ws.Range("H14").Value = "1844.00"
$ws.Range("H15").Value = "2536.65"
$ws.Range("H16").Value = "3200.00"
$ws.Range("H17").Value = "9497.95"
$ws.Range("H18").Value = "1218.00"
$ws.Range("H19").Value = "9961.43"
$ws.Range("H20").Value = "90.00"
$ws.Range("H21").Value = "5654.20"
$ws.Range("H22").Value = "262.31"
$ws.Range("H23").Value = "634.86"
$ws.Range("H24").Value = "160.00"
$ws.Range("H25").Value = "5772.40"
$ws.Range("H26").Value = "10159.43"
$ws.Range("H27").Value = "280.10"
$ws.Range("H28").Value = "120.00"
$ws.Range("H29").Value = "949.00"
$ws.Range("H30").Value = "1647.30"
$ws.Range("H31").Value = "71.38"
$ws.Range("H32").Value = "2874.00"
$ws.Range("H33").Value = "7.22"
$ws.Range("H34").Value = "10732.73"
$ws.Range("H35").Value = "2291.00"
$ws.Range("H36").Value = "21286.69"
$ws.Range("H37").Value = "1960.28"
$ws.Range("H38").Value = "3740.00"
$ws.Range("H39").Value = "9750.00"
$ws.Range("H40").Value = "8100.00"
$ws.Range("H41").Value = "224.19"
$ws.Range("H42").Value = "2790.00"
$ws.Range("H43").Value = "11000.00"
$ws.Range("H44").Value = "433.87"
$ws.Range("H45").Value = "541.56"
$ws.Range("H46").Value = "140.00"
$ws.Range("H47").Value = "34.38"
$ws.Range("H48").Value = "1326.00"
$ws.Range("H49").Value = "339.84"
$ws.Range("H50").Value = "732.58"
$ws.Range("H51").Value = "1195.25"
$ws.Range("H52").Value = "5760.00"
$ws.Range("H53").Value = "25000.00"
$ws.Range("H54").Value = "17028.00"
$ws.Range("H55").Value = "1679.00"
$ws.Range("H56").Value = "1911.02"
$ws.Range("H57").Value = "4890.00"
$ws.Range("H58").Value = "330.00"
$ws.Range("H59").Value = "686.55"
$ws.Range("H60").Value = "3000.00"
$ws.Range("H61").Value = "122721.25"
$ws.Range("H62").Value = "9840.00"
$ws.Range("H63").Value = "45.44"
$ws.Range("H64").Value = "120.00"
$ws.Range("H65").Value = "17.00"
$ws.Range("H66").Value = "171.00"
$ws.Range("H67").Value = "80.00"
$ws.Range("H68").Value = "299.00"
$ws.Range("H69").Value = "19050.38"
$ws.Range("H70").Value = "272.16"
$ws.Range("H71").Value = "4061.20"
$ws.Range("H72").Value = "1742.00"
$ws.Range("H73").Value = "113.17"
$ws.Range("H74").Value = "87.58"
$ws.Range("H75").Value = "16.92"
$ws.Range("H76").Value = "5134.00"
$ws.Range("H77").Value = "25.00"
$ws.Range("H78").Value = "3513.50"
$ws.Range("H79").Value = "14438.00"
$ws.Range("H80").Value = "2267.00"
$ws.Range("H81").Value = "3872.00"
$ws.Range("H82").Value = "90.00"
$ws.Range("H83").Value = "220.00"
$ws.Range("H84").Value = "57262.71"
$ws.Range("H85").Value = "975.00"
$ws.Range("H86").Value = "50.00"
$ws.Range("H87").Value = "1500.00"
$ws.Range("H88").Value = "25000.00"
$ws.Range("H89").Value = "1850.00"
$ws.Range("H90").Value = "3877.70"
$ws.Range("H91").Value = "161.00"
$ws.Range("H92").Value = "9950.00"
$ws.Range("H93").Value = "25.00"
$ws.Range("H94").Value = "268.99"
$ws.Range("H95").Value = "3780.00"
$ws.Range("H96").Value = "2000.00"
$ws.Range("H97").Value = "700.00"
$ws.Range("H98").Value = "1000.00"
$ws.Range("H99").Value = "773.50"
$ws.Range("H100").Value = "250.00"
$ws.Range("H101").Value = "600.00"
$ws.Range("H102").Value = "12000.00"
$ws.Range("H103").Value = "9525.09"
$ws.Range("H104").Value = "1000.00"
$ws.Range("H105").Value = "1900.00"
$ws.Range("H106").Value = "650.00"
$ws.Range("H107").Value = "100.00"
$ws.Range("H108").Value = "2450.00"
$ws.Range("H109").Value = "1500.00"
$ws.Range("H110").Value = "400.00"
$ws.Range("H111").Value = "480.00"
$ws.Range("H112").Value = "4960.00"
$ws.Range("H113").Value = "1000.00"
$ws.Range("H114").Value = "500.00"
$ws.Range("H115").Value = "250.00"
$ws.Range("H116").Value = "105.00"
$ws.Range("H117").Value = "2125.00"
$ws.Range("H118").Value = "200.00"
$ws.Range("H119").Value = "1320.00"
$ws.Range("H120").Value = "1098.50"
$ws.Range("H121").Value = "1520.00"
$ws.Range("H122").Value = "160.24"
$ws.Range("H123").Value = "182.00"
$ws.Range("H124").Value = "50.32"
$ws.Range("H125").Value = "176.00"
$ws.Range("H126").Value = "4280.00"
$ws.Range("H127").Value = "850.00"
$ws.Range("H128").Value = "384.25"
$ws.Range("H129").Value = "35.00"
$ws.Range("H130").Value = "23.64"
$ws.Range("H131").Value = "3120.00"
$ws.Range("H132").Value = "24.60"
$ws.Range("H133").Value = "10725.00"
$ws.Range("H134").Value = "4517.34"
$ws.Range("H135").Value = "1050.00"
$ws.Range("H136").Value = "99.20"
$ws.Range("H137").Value = "357.52"
$ws.Range("H138").Value = "608.00"
$ws.Range("H139").Value = "95.00"
$ws.Range("H140").Value = "289.00"
$ws.Range("H141").Value = "117.72"
$ws.Range("H142").Value = "886.87"
$ws.Range("H143").Value = "3958.90"
$ws.Range("H144").Value = "1174100.82"
$ws.Range("H145").Value = "1149.50"
$ws.Range("H146").Value = "6910.00"
$ws.Range("H147").Value = "2600.00"
$ws.Range("H148").Value = "17000.00"
$ws.Range("H149").Value = "50200.00"
$ws.Range("H150").Value = "10000.00"
$ws.Range("H151").Value = "33300.00"
$ws.Range("H152").Value = "20000.00"
$ws.Range("H153").Value = "187659.00"
$ws.Range("H154").Value = "5970.20"
$ws.Range("H155").Value = "39350.00"
$ws.Range("H156").Value = "15995.00"
$importeRange.ClearFormats()
